$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.315.04"
$ws.Range("E2").Value = "  -0.11%  "

$ws.Range("D3").Value = "2.949.49"
$ws.Range("E3").Value = "  +1.39%  "

$ws.Range("D4").Value = "'0.998"
$ws.Range("E4").Value = "  -0.20%  "

$ws.Range("D5").Value = "'357.22"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'109.84"
$ws.Range("E6").Value = "  -3.75%  "

$ws.Range("D7").Value = "'0.573"
$ws.Range("E7").Value = "  +2.11%  "

$ws.Range("E8").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.628"
$ws.Range("E9").Value = "  +0.05%  "

$ws.Range("D10").Value = "'39.03"
$ws.Range("E10").Value = "  -2.85%  "

$ws.Range("E11").Value = "  +1.48%  "

$ws.Range("D12").Value = "'0.0871"
$ws.Range("E12").Value = "  +0.47%  "

$ws.Range("D13").Value = "'19.56"
$ws.Range("E13").Value = "  -1.63%  "

$ws.Range("D14").Value = "'7.77"
$ws.Range("E14").Value = "  -0.26%  "

$ws.Range("D15").Value = "3.397.26"
$ws.Range("E15").Value = "  +0.83%  "

$ws.Range("D16").Value = "2.946.12"
$ws.Range("E16").Value = "  +0.76%  "

$ws.Range("D17").Value = "'0.980"
$ws.Range("E17").Value = "  -2.27%  "

$ws.Range("D18").Value = "52.214.96"
$ws.Range("E18").Value = "  -0.30%  "

$ws.Range("D19").Value = "'3.52"
$ws.Range("E19").Value = "  +6.16%  "

$ws.Range("D20").Value = "'7.62"
$ws.Range("E20").Value = "  -0.51%  "

$ws.Range("D21").Value = "'13.92"
$ws.Range("E21").Value = "  -1.58%  "

$ws.Range("D22").Value = "0.0₃0984"
$ws.Range("E22").Value = "  +0.47%  "

$ws.Range("D23").Value = "'70.42"
$ws.Range("E23").Value = "  -0.74%  "

$ws.Range("D24").Value = "'270.84"
$ws.Range("E24").Value = "  +0.45%  "

$ws.Range("E25").Value = "  -0.16%  "

$ws.Range("D26").Value = "'0.179"
$ws.Range("E26").Value = "  +2.44%  "

$ws.Range("D27").Value = "'7.87"
$ws.Range("E27").Value = "  +19.36%  "

$ws.Range("D28").Value = "'27.01"
$ws.Range("E28").Value = "  +0.63%  "

$ws.Range("E29").Value = "  +0.10%  "

$ws.Range("E30").Value = "  +6.24%  "

$ws.Range("E31").Value = "  -1.71%  "

$ws.Range("D32").Value = "'37.59"
$ws.Range("E32").Value = "  -1.12%  "

$ws.Range("E33").Value = "  +15.62%  "

$ws.Range("E34").Value = "  -1.62%  "

$ws.Range("D35").Value = "'52.03"
$ws.Range("E35").Value = "  -2.58%  "

$ws.Range("D36").Value = "'0.0443"
$ws.Range("E36").Value = "  -1.99%  "

$ws.Range("E37").Value = "  +0.02%  "

$ws.Range("D38").Value = "'3.20"
$ws.Range("E38").Value = "  -4.88%  "

$ws.Range("D39").Value = "'18.24"
$ws.Range("E39").Value = "  -3.89%  "

$ws.Range("E40").Value = "  -3.77%  "

$ws.Range("E41").Value = "  -1.57%  "

$ws.Range("D42").Value = "'0.121"
$ws.Range("E42").Value = "  +2.62%  "

$ws.Range("D43").Value = "'22.84"
$ws.Range("E43").Value = "  -1.74%  "

$ws.Range("D44").Value = "'119.37"
$ws.Range("E44").Value = "  -0.91%  "

$ws.Range("E45").Value = "  -1.17%  "

$ws.Range("D48").Value = "2.137.17"
$ws.Range("E48").Value = "  -2.24%  "

$ws.Range("E49").Value = "  -5.91%  "

$ws.Range("E50").Value = "  +0.64%  "

$ws.Range("D51").Value = "'0.913"
$ws.Range("E51").Value = "  -4.54%  "

$ws.Range("B46").Value = "ApeXProtocol"
$ws.Range("C46").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D46").Value = "'2.47"
$ws.Range("E46").Value = "  -4.98%  "

$ws.Range("B47").Value = "NEARProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D47").Value = "'3.45"
$ws.Range("E47").Value = "  -2.36%  "
